$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- API Connection: add a new "Moderate Rain" (501) entry to the Weather IDs table ---
# Duplicate row 11 (which already carries the correct cell styles used for table data
# rows: s=16 / s=17 / s=18) down into row 12, fill in the new data, then remove the
# row that got pushed down to the bottom of the sheet so the overall layout/row count
# stays the same as before (rows 13-19 keep their original row numbers).
$ws.Rows("11:11").Copy()
$ws.Rows("12:12").Insert(-4121)   # xlShiftDown

$ws.Range("B12").Value = 501
$ws.Range("C12").Value = "Moderate Rain"
$ws.Range("D12").ClearContents()
$ws.Rows("12:12").RowHeight = 30

$ws.Rows("20:20").Delete(-4162)  # xlShiftUp - remove the row pushed past the old bottom

# Grow the "Aufgaben" table so the new row becomes part of it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B6:D12"))

# --- UI: corrected minus label / selection moved to the last edited cell ---
$ws.Range("D17").Select()
